$wb = $excel.ActiveWorkbook

# --- MAIN_CONTROLLER sheet -------------------------------------------------
# A stray "N" (same RunStatus flag used in column B) is recorded in H14,
# which was already the (multi-)selected cell before the edit.
$wsMain = $wb.Worksheets.Item("MAIN_CONTROLLER")
$wsMain.Cells.Item(14, 8).Value = "N"

# --- MOBILE_CONFIGURATION sheet --------------------------------------------
# Package name switched from the QA build to the pre-prod build for the
# "scroll up/down till element" rows (3 & 4).
$wsMobile = $wb.Worksheets.Item("MOBILE_CONFIGURATION")
$wsMobile.Cells.Item(3, 4).Value = "com.mahindra.fospreprod"
$wsMobile.Cells.Item(4, 4).Value = "com.mahindra.fospreprod"
# Column D widens to fit the new (longer) value.
$wsMobile.Columns.Item(4).ColumnWidth = 23.5

# --- DATASHEET sheet ---------------------------------------------------------
# Test data workbook reference updated to the scroll-test data file.
$wsData = $wb.Worksheets.Item("DATASHEET")
$wsData.Cells.Item(3, 4).Value = "FOSScroll.xlsx"

# --- Restore on-screen selections for each sheet and leave DATASHEET as the
#     active tab (selecting a sheet last makes it the active one on save).
$wsMain.Range("E2").Select()
$wsMobile.Range("D10").Select()
$wsData.Range("D13").Select()

Write-Host "done"
